$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily log rows continuing the series (dates 44256-44263),
# matching the existing table layout: A=date, B:G=exercise reps.
$rows = @(
    @(44256, 0,  0,  0,  0,   0,   0),
    @(44257, 0,  0,  0,  0,   0,   0),
    @(44258, 24, 30, 0,  0,   0,   0),
    @(44259, 0,  0,  0,  0,   0,   24),
    @(44260, 0,  0,  0,  0,   0,   0),
    @(44261, 0,  0,  0,  0,   0,   0),
    @(44262, 24, 0,  0,  0,   0,   0),
    @(44263, 20, 30, 0,  30,  100, 0)
)

$startRow = 42
$lastExistingRow = 41

# Grab the date formatting already used in column A so the new rows
# match the existing table's style (same <c s="1"> as the rows above).
$ws.Cells.Item($lastExistingRow, 1).Copy() | Out-Null

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $dateCell.Value2 = $data[0]

    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $data[$c]
    }
}
$excel.CutCopyMode = $false

# Match the new selection/view state recorded for the edit.
$ws.Range("E49").Select()
